# Insert a new row before row 39, duplicating the phone number (as a true
# number this time, without the leading zero) and birthday of the
# existing "09876543" customer, with 0 points. The original row (text
# phone "09876543") shifts down to row 40 and loses its birthday value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 39 (and everything below) down to make room for the new row.
$ws.Rows.Item(39).Insert()

# New row 39: numeric phone, same birthday, 0 points.
$ws.Cells.Item(39, 1).Value = 9876543
$ws.Cells.Item(39, 2).NumberFormat = "@"
$ws.Cells.Item(39, 2).Value = "2025-08-12"
$ws.Cells.Item(39, 3).Value = 0

# Row 40 is the original row, now shifted down; clear its birthday.
$ws.Cells.Item(40, 2).Value = ""
